# Fix species labeling error in carter_2005_temp_thresholds.xlsx
# Previously every row was labeled "Chinook"; the second half of the
# data actually corresponds to steelhead. Also normalize casing of
# the Chinook label to lowercase "chinook" to match "steelhead".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-41: species = chinook (lowercase)
$ws.Range("A2:A41").Value = "chinook"

# Rows 42-71: species = steelhead
$ws.Range("A42:A71").Value = "steelhead"

# Leave the selection on the newly-corrected steelhead block, matching
# where the editor was working when the fix was made.
$ws.Range("A42:A71").Select()
